$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 0.453125
    3 = 0.40625
    4 = 0.34375
    5 = 0.328125
    6 = 0.328125
    7 = 0.328125
    8 = 0.3125
    9 = 0.296875
    10 = 0.28125
    11 = 0.265625
    12 = 0.328125
    13 = 0.3125
    14 = 0.28125
    15 = 0.359375
    16 = 0.328125
    17 = 0.265625
    18 = 0.25
    19 = 0.28125
    20 = 0.3125
    21 = 0.28125
    22 = 0.328125
    23 = 0.3125
    24 = 0.3125
    25 = 0.3125
    26 = 0.296875
    27 = 0.296875
    28 = 0.296875
    29 = 0.296875
    30 = 0.296875
    31 = 0.296875
    32 = 0.296875
    33 = 0.3125
    34 = 0.3125
    35 = 0.3125
    36 = 0.296875
    37 = 0.296875
    38 = 0.296875
    39 = 0.296875
    40 = 0.296875
    41 = 0.28125
    42 = 0.28125
    43 = 0.28125
    44 = 0.28125
    45 = 0.28125
    46 = 0.28125
    47 = 0.28125
    48 = 0.28125
    49 = 0.28125
    50 = 0.28125
    51 = 0.28125
    52 = 0.28125
    53 = 0.28125
    54 = 0.28125
    55 = 0.28125
    56 = 0.28125
    57 = 0.28125
    58 = 0.28125
    59 = 0.28125
    60 = 0.28125
    61 = 0.28125
    62 = 0.28125
    63 = 0.28125
    64 = 0.28125
    65 = 0.28125
    66 = 0.28125
    67 = 0.28125
    68 = 0.28125
    69 = 0.28125
    70 = 0.28125
    71 = 0.28125
    72 = 0.28125
    73 = 0.28125
    74 = 0.28125
    75 = 0.28125
    76 = 0.28125
    77 = 0.28125
    78 = 0.28125
    79 = 0.28125
    80 = 0.28125
    81 = 0.28125
    82 = 0.28125
    83 = 0.28125
    84 = 0.28125
    85 = 0.28125
    86 = 0.28125
    87 = 0.28125
    88 = 0.28125
    89 = 0.28125
    90 = 0.28125
    91 = 0.28125
    92 = 0.28125
    93 = 0.28125
    94 = 0.28125
    95 = 0.28125
    96 = 0.28125
    97 = 0.28125
    98 = 0.28125
    99 = 0.28125
    100 = 0.28125
    101 = 0.28125
    102 = 0.28125
    103 = 0.25
    104 = 0.171875
    105 = 0.234375
    106 = 0.21875
    107 = 0.140625
    108 = 0.28125
    109 = 0.21875
    110 = 0.171875
    111 = 0.234375
    112 = 0.203125
    113 = 0.234375
    114 = 0.109375
    115 = 0.15625
    116 = 0.125
    117 = 0.203125
    118 = 0.1967213114754098
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 2).Value = $values[$r]
}

$ws.Range("A2").Select()
$ws.Range("A1:XFD1048576").Select()